$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header values in row 2 and add a new column D
$ws.Range("A2").Value = "BusinessKey"
$ws.Range("B2").Value = "Code"
$ws.Range("C2").Value = "CustomReportType_ID"
$ws.Range("D2").Value = "Name"

# Match style of existing header cells for new D2 cell
$ws.Range("D2").Style = $ws.Range("C2").Style
